# Refresh the cryptocurrency Price (column D) and 1h Volume-change (column E)
# figures to match the latest coinranking.com scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.958.66'
$ws.Range("E2").Value = '  +1.76%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.153.90'
$ws.Range("E3").Value = '  +3.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.14'
$ws.Range("E5").Value = '  +2.64%  '

$ws.Range("E6").Value = '  +5.18%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.150.95'
$ws.Range("E8").Value = '  +3.08%  '

$ws.Range("E9").Value = '  +1.82%  '

$ws.Range("E10").Value = '  +4.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.10'
$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("E12").Value = '  +3.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000264'
$ws.Range("E13").Value = '  +14.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.04'
$ws.Range("E14").Value = '  +4.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.665.88'
$ws.Range("E15").Value = '  +2.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.986.10'
$ws.Range("E16").Value = '  +1.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.159.91'
$ws.Range("E17").Value = '  +3.32%  '

$ws.Range("E18").Value = '  +4.57%  '

$ws.Range("E19").Value = '  +1.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '506.01'
$ws.Range("E20").Value = '  +3.98%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.83'
$ws.Range("E21").Value = '  +3.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.717'
$ws.Range("E22").Value = '  +3.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.24'
$ws.Range("E23").Value = '  +3.69%  '

$ws.Range("E24").Value = '  +2.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.19'
$ws.Range("E25").Value = '  +2.06%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  +3.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.86'
$ws.Range("E28").Value = '  +8.50%  '

$ws.Range("E29").Value = '  +5.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.80'
$ws.Range("E30").Value = '  +9.67%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.56'
$ws.Range("E31").Value = '  +4.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.06%  '

$ws.Range("E33").Value = '  +2.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.17'
$ws.Range("E34").Value = '  +7.76%  '

$ws.Range("E35").Value = '  +3.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '54.69'
$ws.Range("E36").Value = '  -0.46%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0898'
$ws.Range("E37").Value = '  +9.96%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '463.30'
$ws.Range("E38").Value = '  +4.84%  '

$ws.Range("E39").Value = '  +1.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.98'
$ws.Range("E40").Value = '  +9.05%  '

$ws.Range("E41").Value = '  +3.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.044.36'
$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("E43").Value = '  +0.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.45'
$ws.Range("E44").Value = '  +9.27%  '

$ws.Range("E45").Value = '  +2.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.49'
$ws.Range("E46").Value = '  +2.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0585'
$ws.Range("E47").Value = '  +13.14%  '

$ws.Range("E48").Value = '  -0.04%  '

$ws.Range("E49").Value = '  +0.51%  '

$ws.Range("E50").Value = '  +4.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.24'
$ws.Range("E51").Value = '  +1.53%  '
